$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update Hydrogen / Iron & steel value, clear Non-metallic minerals value
$ws.Range("B3").Value = 810206.2497117716
$ws.Range("D3").Value = $null

# Row 4: update Methanol / Chemicals value
$ws.Range("C4").Value = 45.55891441616917

# Row 5: update Ammonia / Chemicals value
$ws.Range("C5").Value = 2860.383513718285

# Row 7: rename "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 691.7720860150997

# New row 8: "Other" row, re-created below Biogas
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 526.9099049891719

# Copy the row-label formatting (bold, borders, centered) from A7 onto the new A8 cell
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
